$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 79 was previously the last row (special "latest" date format);
# now that a new row is appended, row 79 reverts to the normal
# date-time number format used by all the other non-final rows.
$ws.Range("A79").NumberFormat = $ws.Range("A78").NumberFormat

# Append the new day's data (row 80), keeping the "latest row" date
# format that row 79 used to have.
$ws.Range("A80").Value = 45820
$ws.Range("A80").NumberFormat = "YYYY-MM-DD"

$ws.Range("B80").Value = 339
$ws.Range("C80").Value = 341
$ws.Range("D80").Value = 345
